# Append: 2025-10-30 12:49 JST
# Update the acquisition-timestamp column (A) for rows 2-11 on the first
# worksheet ("ランサーズ") from the old run time to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "2025-10-30 12:37:03") {
        $cell.Value2 = "2025-10-30 12:49:01"
    }
}
